$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1997.625
$ws.Range("J17").Value = 1754.463
$ws.Range("L17").Value = 5263.389
$ws.Range("N17").Value = -5599.389
$ws.Range("H43").Value = 1394.625
$ws.Range("I43").Value = 1500
$ws.Range("J43").Value = 1379.5714
$ws.Range("K43").Value = 1500
$ws.Range("L43").Value = 1379.5714
$ws.Range("M43").Value = -1431
$ws.Range("N43").Value = -1517.5714
$ws.Range("H69").Value = 6750
$ws.Range("J69").Value = 8000
$ws.Range("L69").Value = 24000
$ws.Range("N69").Value = -25748
$ws.Range("H72").Value = 6750
$ws.Range("J72").Value = 8000
$ws.Range("L72").Value = 72000
$ws.Range("N72").Value = -80736
$ws.Range("H101").Value = 1779843.4
$ws.Range("I101").Value = 4444609
$ws.Range("K101").Value = 13333827
$ws.Range("M101").Value = -13332205
$ws.Range("H132").Value = 1144.9846
$ws.Range("I132").Value = 897.7037
$ws.Range("K132").Value = 2693.1111
$ws.Range("M132").Value = -163.1111000000001
$ws.Range("H135").Value = 362.875
$ws.Range("I135").Value = 321.41025
$ws.Range("J135").Value = 1980
$ws.Range("K135").Value = 2892.69225
$ws.Range("L135").Value = 17820
$ws.Range("M135").Value = -357.6922500000001
$ws.Range("N135").Value = -22890
$ws.Range("H137").Value = 25562.414
$ws.Range("I137").Value = 975.6923
$ws.Range("J137").Value = 68179.39999999999
$ws.Range("K137").Value = 2927.0769
$ws.Range("L137").Value = 204538.2
$ws.Range("M137").Value = -377.0769
$ws.Range("N137").Value = -209638.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3516.962
$ws.Range("I32").Value = 3063.5657
$ws.Range("J32").Value = 15003
$ws.Range("K32").Value = 3063.5657
$ws.Range("L32").Value = 15003
$ws.Range("M32").Value = -2776.5657
$ws.Range("N32").Value = -15577
$ws.Range("H61").Value = 3696.3333
$ws.Range("I61").Value = 1129.1904
$ws.Range("J61").Value = 21666.334
$ws.Range("K61").Value = 1129.1904
$ws.Range("L61").Value = 21666.334
$ws.Range("M61").Value = -917.1904
$ws.Range("N61").Value = -22090.334
$ws.Range("H74").Value = 1382.6957
$ws.Range("I74").Value = 955.8889
$ws.Range("J74").Value = 2919.2
$ws.Range("K74").Value = 955.8889
$ws.Range("L74").Value = 2919.2
$ws.Range("M74").Value = -81.88890000000004
$ws.Range("N74").Value = -4667.2
$ws.Range("H77").Value = 1382.6957
$ws.Range("I77").Value = 955.8889
$ws.Range("J77").Value = 2919.2
$ws.Range("K77").Value = 4779.444500000001
$ws.Range("L77").Value = 14596
$ws.Range("M77").Value = -411.4445000000005
$ws.Range("N77").Value = -23332
$ws.Range("H88").Value = 2887.4546
$ws.Range("I88").Value = 1555
$ws.Range("J88").Value = 3020.7
$ws.Range("K88").Value = 1555
$ws.Range("L88").Value = 3020.7
$ws.Range("M88").Value = -1149
$ws.Range("N88").Value = -3832.7
$ws.Range("H91").Value = 2887.4546
$ws.Range("I91").Value = 1555
$ws.Range("J91").Value = 3020.7
$ws.Range("K91").Value = 1555
$ws.Range("L91").Value = 3020.7
$ws.Range("M91").Value = -151
$ws.Range("N91").Value = -5828.7
$ws.Range("H97").Value = 880.6
$ws.Range("I97").Value = 845.1053000000001
$ws.Range("K97").Value = 845.1053000000001
$ws.Range("M97").Value = -349.1053000000001
$ws.Range("H122").Value = 1672.4286
$ws.Range("I122").Value = 1474.5
$ws.Range("J122").Value = 2167.25
$ws.Range("K122").Value = 4423.5
$ws.Range("L122").Value = 6501.75
$ws.Range("M122").Value = -1973.5
$ws.Range("N122").Value = -11401.75
$ws.Range("H132").Value = 1851.674
$ws.Range("I132").Value = 1393.7632
$ws.Range("K132").Value = 4181.2896
$ws.Range("M132").Value = -1651.2896
$ws.Range("H136").Value = 3696.3333
$ws.Range("I136").Value = 1129.1904
$ws.Range("J136").Value = 21666.334
$ws.Range("K136").Value = 3387.5712
$ws.Range("L136").Value = 64999.00199999999
$ws.Range("M136").Value = -837.5711999999999
$ws.Range("N136").Value = -70099.00199999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1891.9667
$ws.Range("I20").Value = 1528.8695
$ws.Range("J20").Value = 3085
$ws.Range("K20").Value = 1528.8695
$ws.Range("L20").Value = 3085
$ws.Range("M20").Value = -1281.8695
$ws.Range("N20").Value = -3579
$ws.Range("H94").Value = 1034
$ws.Range("I94").Value = 1480.6
$ws.Range("K94").Value = 1480.6
$ws.Range("M94").Value = -1029.6
$ws.Range("H134").Value = 6897.5454
$ws.Range("I134").Value = 7879.1763
$ws.Range("K134").Value = 23637.5289
$ws.Range("M134").Value = -21102.5289
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1476.4894
$ws.Range("I31").Value = 1048.7179
$ws.Range("J31").Value = 3561.875
$ws.Range("K31").Value = 1048.7179
$ws.Range("L31").Value = 3561.875
$ws.Range("M31").Value = -753.7179000000001
$ws.Range("N31").Value = -4151.875
$ws.Range("H34").Value = 1476.4894
$ws.Range("I34").Value = 1048.7179
$ws.Range("J34").Value = 3561.875
$ws.Range("K34").Value = 1048.7179
$ws.Range("L34").Value = 3561.875
$ws.Range("M34").Value = -846.7179000000001
$ws.Range("N34").Value = -3965.875
$ws.Range("H58").Value = 750535.1
$ws.Range("I58").Value = 1061085.5
$ws.Range("K58").Value = 1061085.5
$ws.Range("M58").Value = -1060882.5
$ws.Range("H92").Value = 44999
$ws.Range("J92").Value = 44999
$ws.Range("L92").Value = 44999
$ws.Range("N92").Value = -49991
$ws.Range("H132").Value = 1603.6111
$ws.Range("I132").Value = 1103.9
$ws.Range("J132").Value = 3031.3572
$ws.Range("K132").Value = 3311.7
$ws.Range("L132").Value = 9094.071599999999
$ws.Range("M132").Value = -781.7000000000003
$ws.Range("N132").Value = -14154.0716
$ws.Range("H134").Value = 1459.7142
$ws.Range("I134").Value = 1377.3695
$ws.Range("J134").Value = 1682.5294
$ws.Range("K134").Value = 4132.1085
$ws.Range("L134").Value = 5047.5882
$ws.Range("M134").Value = -1597.1085
$ws.Range("N134").Value = -10117.5882
$ws.Range("H136").Value = 750535.1
$ws.Range("I136").Value = 1061085.5
$ws.Range("K136").Value = 3183256.5
$ws.Range("M136").Value = -3180706.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 6492.4546
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null
$ws.Range("H140").Value = 2122.75
$ws.Range("J140").Value = 2697.3333
$ws.Range("L140").Value = 8091.999899999999
$ws.Range("N140").Value = -18451.9999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1067.9584
$ws.Range("I97").Value = 1203.6428
$ws.Range("K97").Value = 1203.6428
$ws.Range("M97").Value = -707.6428000000001
$ws.Range("H132").Value = 940480.9
$ws.Range("J132").Value = 2936.353
$ws.Range("L132").Value = 8809.059000000001
$ws.Range("N132").Value = -13869.059
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 15000
$ws.Range("J20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("N20").Value = -15452
$ws.Range("H93").Value = 1063.4286
$ws.Range("I93").Value = 820
$ws.Range("K93").Value = 820
$ws.Range("M93").Value = 428
$ws.Range("H132").Value = 3648.4167
$ws.Range("I132").Value = 2321.75
$ws.Range("J132").Value = 6301.75
$ws.Range("K132").Value = 6965.25
$ws.Range("L132").Value = 18905.25
$ws.Range("M132").Value = -4435.25
$ws.Range("N132").Value = -23965.25
$ws.Range("H136").Value = 1893.1538
$ws.Range("J136").Value = 5498
$ws.Range("L136").Value = 16494
$ws.Range("N136").Value = -21594
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 9100
$ws.Range("I21").Value = 2750
$ws.Range("J21").Value = 13333.333
$ws.Range("K21").Value = 2750
$ws.Range("L21").Value = 13333.333
$ws.Range("M21").Value = -2515
$ws.Range("N21").Value = -13803.333
$ws.Range("H28").Value = 10000
$ws.Range("J28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("N28").Value = -10696
$ws.Range("H35").Value = 9100
$ws.Range("I35").Value = 2750
$ws.Range("J35").Value = 13333.333
$ws.Range("K35").Value = 2750
$ws.Range("L35").Value = 13333.333
$ws.Range("M35").Value = -2460
$ws.Range("N35").Value = -13913.333
$ws.Range("H41").Value = 13990.75
$ws.Range("J41").Value = 13990.75
$ws.Range("L41").Value = 13990.75
$ws.Range("N41").Value = -14770.75
$ws.Range("H122").Value = 51735
$ws.Range("I122").Value = 69192.87
$ws.Range("J122").Value = 1543.625
$ws.Range("K122").Value = 207578.61
$ws.Range("L122").Value = 4630.875
$ws.Range("M122").Value = -205128.61
$ws.Range("N122").Value = -9530.875
